$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Date Apply" (column D) values for rows 5-7 and 8-10
$ws.Range("D5").Value = "2023-01-27 02:01:03"
$ws.Range("D6").Value = "2023-01-27 02:01:03"
$ws.Range("D7").Value = "2023-01-27 02:01:03"
$ws.Range("D8").Value = "2023-13-29 02:13:51"
$ws.Range("D9").Value = "2023-13-29 02:13:51"
$ws.Range("D10").Value = "2023-13-29 02:13:51"

# "Status" (column E) values for rows 5-10
$ws.Range("E5").Value = "No Responce"
$ws.Range("E6").Value = "No Responce"
$ws.Range("E7").Value = "No Responce"
$ws.Range("E8").Value = "No Responce"
$ws.Range("E9").Value = "No Responce"
$ws.Range("E10").Value = "No Responce"
